# Add a new row (row 11) to Sheet1 with the "length of mooring chain" (锚链长度)
# info, matching the style used by the other label rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous label cell (A10) onto the new label
# cell (A11), so it picks up the same left-aligned style used by the other
# row labels in column A. The new data cells (B11:D11) keep the default
# (unstyled) formatting, matching the rest of the data cells added.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A11").Value = "锚链长度"
$ws.Range("B11").Value = 15.2367
$ws.Range("C11").Value = 21.7229
$ws.Range("D11").Value = 22.05

# Match the final selection recorded in the saved workbook.
$ws.Range("G11").Select()
